# Commit: "updated files with strain names"
#
# The harvester for every data row (rows 2-43 of Sheet1) is changed from
# "S.GISH" to "H.BROWN". All other per-row data (experimentDesign, strain,
# genotype, treatment, timePoint, replicate, ...) is left untouched - the
# shared-strings index churn visible in the raw XML diff is just an
# automatic side effect of "S.GISH" dropping out of the shared strings
# table and "H.BROWN" being appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: the source workbook also turns on iterative calculation
# (calcPr iterateDelta="1E-4"). Set it via the Application object in case
# the host wires it through; harmless no-op otherwise.
try {
    $excel.Iteration = $true
    $excel.MaxChange = 0.0001
} catch {
}

for ($row = 2; $row -le 43; $row++) {
    $ws.Range("B$row").Value = "H.BROWN"
}

# Match the author's final selection (B3:B43, active cell B3).
$ws.Range("B3:B43").Select()
